# "only fixed get price" - append two new log rows produced by the
# get_price command to the bottom of the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42 ---------------------------------------------------------
$ws.Range("A42").Value = "2024-10-09 14:10:23"
$ws.Range("B42").Value = "get_price"
$ws.Range("C42").Value = "https://example.com/product"
# Leading apostrophe forces these to be stored as text (not auto-converted
# to a number / date) while keeping the underlying display value intact.
$ws.Range("D42").Value = "'100.00"
$ws.Range("E42").Value = "'2024-10-09"
$ws.Range("F42").Value = "14:10:23"

# --- Row 43 ---------------------------------------------------------
$ws.Range("A43").Value = "2024-10-09 14:22:04"
$ws.Range("B43").Value = "get_price"
$ws.Range("C43").Value = "https://example.com/product"
$ws.Range("D43").Value = "'100.00"
$ws.Range("E43").Value = "'2024-10-09"
$ws.Range("F43").Value = "14:22:04"

# The quote-prefix entry above leaves a "quotePrefix" style behind; reset
# those cells back to the plain/default style so they match the rest of
# the table (which carries no explicit cell style).
$ws.Range("D42:E43").Style = "Normal"
